$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price (D) and Volume(1h) (E) columns store values as text,
# so force a Text number format before writing each string value
# to prevent Excel from auto-converting it to a number/percentage.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "301.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.44%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.14"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.75%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.997"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.04%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07851"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.66%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.116"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "-9.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.833"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "1.27%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.833"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.89%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9250"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "0.08%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1757"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.72%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07750"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.26%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08624"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-3.16%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03157"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.89%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001520"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.60%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005820"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-5.56%"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2,107.56%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.477"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-0.27%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "2.117"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-7.58%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3277"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.15%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1319"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "-0.65%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.265"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2.88%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.1855"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "9.14%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.04600"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.28%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001223"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.95%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004447"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.11%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001246"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.81%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01758"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "1.14%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04805"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.03%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007525"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "8.36%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1360"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.37%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002254"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "2.91%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.74%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006021"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-4.13%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00000000747"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.34%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-29.29%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002093"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.34%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0001993"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.34%"
